$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '52.282.07'
$ws.Range('E2').Value = '  +1.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.847.42'
$ws.Range('E3').Value = '  +3.51%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '362.85'
$ws.Range('E5').Value = '  +9.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '115.97'
$ws.Range('E6').Value = '  -0.92%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.553'
$ws.Range('E7').Value = '  +3.72%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  +5.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.47'
$ws.Range('E10').Value = '  +2.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0865'
$ws.Range('E11').Value = '  +3.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.17'
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('E13').Value = '  +1.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.89'
$ws.Range('E14').Value = '  +3.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.298.93'
$ws.Range('E15').Value = '  +3.90%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.870.88'
$ws.Range('E16').Value = '  +4.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.905'
$ws.Range('E17').Value = '  +2.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '52.283.52'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.95'
$ws.Range('E19').Value = '  +3.02%  '
$ws.Range('B20').Value = 'ImmutableX'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.19'
$ws.Range('E20').Value = '  +5.31%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.32'
$ws.Range('E21').Value = '  +7.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0000100'
$ws.Range('E22').Value = '  +4.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.56'
$ws.Range('E23').Value = '  +1.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '271.61'
$ws.Range('E24').Value = '  -2.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.84'
$ws.Range('E25').Value = '  +7.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.24'
$ws.Range('E26').Value = '  +1.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.32'
$ws.Range('E28').Value = '  +0.98%  '
$ws.Range('E29').Value = '  +1.44%  '
$ws.Range('E30').Value = '  +0.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.61'
$ws.Range('E31').Value = '  -0.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '51.28'
$ws.Range('E32').Value = '  +1.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.87'
$ws.Range('E33').Value = '  +5.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0442'
$ws.Range('E34').Value = '  +27.34%  '
$ws.Range('E35').Value = '  +1.95%  '
$ws.Range('E36').Value = '  +2.50%  '
$ws.Range('E37').Value = '  +0.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.01'
$ws.Range('E38').Value = '  +1.49%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.29'
$ws.Range('E39').Value = '  +2.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.78'
$ws.Range('E40').Value = '  -1.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.61'
$ws.Range('E41').Value = '  +9.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '23.62'
$ws.Range('E42').Value = '  +1.96%  '
$ws.Range('E43').Value = '  +2.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '127.38'
$ws.Range('E44').Value = '  -2.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.29'
$ws.Range('E45').Value = '  +1.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.39'
$ws.Range('E46').Value = '  +1.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.076.99'
$ws.Range('E47').Value = '  -1.46%  '
$ws.Range('E48').Value = '  +2.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.954'
$ws.Range('E49').Value = '  +8.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.63'
$ws.Range('E50').Value = '  +1.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.03'
$ws.Range('E51').Value = '  +0.89%  '
